$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 824.1177
$ws.Range("I92").Value = 571.9286
$ws.Range("J92").Value = 2001
$ws.Range("K92").Value = 571.9286
$ws.Range("L92").Value = 2001
$ws.Range("M92").Value = 676.0714
$ws.Range("N92").Value = -4497

$ws.Range("H104").Value = 1055.3334
$ws.Range("I104").Value = 1055.3334
$ws.Range("K104").Value = 3166.0002
$ws.Range("M104").Value = -1419.0002

$ws.Range("H137").Value = 1888614.9
$ws.Range("I137").Value = 3031626.5
$ws.Range("J137").Value = 2645.9
$ws.Range("K137").Value = 9094879.5
$ws.Range("L137").Value = 7937.700000000001
$ws.Range("M137").Value = -9092329.5
$ws.Range("N137").Value = -13037.7

$ws.Range("H138").Value = 4110553.2
$ws.Range("I138").Value = 386114.03
$ws.Range("J138").Value = 9806754
$ws.Range("K138").Value = 1158342.09
$ws.Range("L138").Value = 29420262
$ws.Range("M138").Value = -1153202.09
$ws.Range("N138").Value = -29430542

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H21").Value = 5590.5557
$ws.Range("I21").Value = 3885.8333
$ws.Range("J21").Value = 9000
$ws.Range("K21").Value = 3885.8333
$ws.Range("L21").Value = 9000
$ws.Range("M21").Value = -3511.8333
$ws.Range("N21").Value = -9748

$ws.Range("H23").Value = 30006
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 30006
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 30006
$ws.Range("M23").ClearContents()
$ws.Range("N23").Value = -30524

$ws.Range("H45").Value = 1205.5264
$ws.Range("I45").Value = 892.3333
$ws.Range("J45").Value = 1742.4286
$ws.Range("K45").Value = 892.3333
$ws.Range("L45").Value = 1742.4286
$ws.Range("M45").Value = -515.3333
$ws.Range("N45").Value = -2496.4286

$ws.Range("H63").Value = 5984.1665
$ws.Range("I63").Value = 3181
$ws.Range("K63").Value = 3181
$ws.Range("M63").Value = -2495

$ws.Range("H66").Value = 5984.1665
$ws.Range("I66").Value = 3181
$ws.Range("K66").Value = 15905
$ws.Range("M66").Value = -12473

$ws.Range("H74").Value = 4501572.5
$ws.Range("I74").Value = 5977268
$ws.Range("J74").Value = 74485.86
$ws.Range("K74").Value = 5977268
$ws.Range("L74").Value = 74485.86
$ws.Range("M74").Value = -5976394
$ws.Range("N74").Value = -76233.86

$ws.Range("H77").Value = 4501572.5
$ws.Range("I77").Value = 5977268
$ws.Range("J77").Value = 74485.86
$ws.Range("K77").Value = 29886340
$ws.Range("L77").Value = 372429.3
$ws.Range("M77").Value = -29881972
$ws.Range("N77").Value = -381165.3

$ws.Range("H97").Value = 1954134
$ws.Range("I97").Value = 2501077.5
$ws.Range("K97").Value = 2501077.5
$ws.Range("M97").Value = -2500581.5

$ws.Range("H102").Value = 7144575.5
$ws.Range("I102").Value = 8404853
$ws.Range("K102").Value = 8404853
$ws.Range("M102").Value = -8403231

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H12").Value = 1000
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 1000
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 1000
$ws.Range("M12").ClearContents()
$ws.Range("N12").Value = -1336

$ws.Range("H24").Value = 1413.6666
$ws.Range("I24").Value = 1413.6666
$ws.Range("J24").Value = 0
$ws.Range("K24").Value = 1413.6666
$ws.Range("L24").Value = 0
$ws.Range("M24").Value = -1178.6666
$ws.Range("N24").ClearContents()

$ws.Range("H99").Value = 1367
$ws.Range("I99").Value = 1438.2354
$ws.Range("J99").Value = 1165.1666
$ws.Range("K99").Value = 1438.2354
$ws.Range("L99").Value = 1165.1666
$ws.Range("M99").Value = 59.76459999999997
$ws.Range("N99").Value = -4161.1666

$ws.Range("H134").Value = 2146.6875
$ws.Range("I134").Value = 1303
$ws.Range("K134").Value = 3909
$ws.Range("M134").Value = -1374

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H19").Value = 731.1429000000001
$ws.Range("I19").Value = 454.75
$ws.Range("J19").Value = 1099.6666
$ws.Range("K19").Value = 454.75
$ws.Range("L19").Value = 1099.6666
$ws.Range("M19").Value = -284.75
$ws.Range("N19").Value = -1439.6666

$ws.Range("H24").Value = 731.1429000000001
$ws.Range("I24").Value = 454.75
$ws.Range("J24").Value = 1099.6666
$ws.Range("K24").Value = 454.75
$ws.Range("L24").Value = 1099.6666
$ws.Range("M24").Value = -284.75
$ws.Range("N24").Value = -1439.6666

$ws.Range("H32").Value = 3000
$ws.Range("I32").Value = 3000
$ws.Range("K32").Value = 3000
$ws.Range("M32").Value = -2684

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 1126.125
$ws.Range("I34").Value = 834.3333
$ws.Range("J34").Value = 2001.5
$ws.Range("K34").Value = 2502.9999
$ws.Range("L34").Value = 6004.5
$ws.Range("M34").Value = -2418.9999
$ws.Range("N34").Value = -6172.5

$ws.Range("H96").Value = 5012.5
$ws.Range("I96").Value = 0
$ws.Range("J96").Value = 5012.5
$ws.Range("K96").Value = 0
$ws.Range("L96").Value = 15037.5
$ws.Range("M96").ClearContents()
$ws.Range("N96").Value = -19155.5

$ws.Range("H131").Value = 11628913
$ws.Range("I131").Value = 71428936
$ws.Range("J131").Value = 1130.6666
$ws.Range("K131").Value = 214286808
$ws.Range("L131").Value = 3391.9998
$ws.Range("M131").Value = -214281768
$ws.Range("N131").Value = -13471.9998

$ws.Range("H132").Value = 879.5714
$ws.Range("I132").Value = 650.8
$ws.Range("J132").Value = 1006.6667
$ws.Range("K132").Value = 5857.2
$ws.Range("L132").Value = 9060.0003
$ws.Range("M132").Value = -3327.2
$ws.Range("N132").Value = -14120.0003

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1982
$ws.Range("I97").Value = 2023.125
$ws.Range("J97").Value = 1817.5
$ws.Range("K97").Value = 2023.125
$ws.Range("L97").Value = 1817.5
$ws.Range("M97").Value = -1527.125
$ws.Range("N97").Value = -2809.5

$ws.Range("H122").Value = 3030.9033
$ws.Range("I122").Value = 2906.4546
$ws.Range("J122").Value = 3335.111
$ws.Range("K122").Value = 8719.363799999999
$ws.Range("L122").Value = 10005.333
$ws.Range("M122").Value = -6269.363799999999
$ws.Range("N122").Value = -14905.333

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1019.6923
$ws.Range("I22").Value = 694.44446
$ws.Range("J22").Value = 1751.5
$ws.Range("K22").Value = 694.44446
$ws.Range("L22").Value = 1751.5
$ws.Range("M22").Value = -399.44446
$ws.Range("N22").Value = -2341.5

$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("M23").ClearContents()

$ws.Range("H27").Value = 1019.6923
$ws.Range("I27").Value = 694.44446
$ws.Range("J27").Value = 1751.5
$ws.Range("K27").Value = 694.44446
$ws.Range("L27").Value = 1751.5
$ws.Range("M27").Value = -587.44446
$ws.Range("N27").Value = -1965.5

$ws.Range("H74").Value = 5020098.5
$ws.Range("J74").Value = 40000
$ws.Range("L74").Value = 40000
$ws.Range("N74").Value = -41996

$ws.Range("H77").Value = 5020098.5
$ws.Range("J77").Value = 40000
$ws.Range("L77").Value = 120000
$ws.Range("N77").Value = -129984

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1268.75
$ws.Range("I126").Value = 1010
$ws.Range("J126").Value = 1700
$ws.Range("K126").Value = 3030
$ws.Range("L126").Value = 5100
$ws.Range("M126").Value = -560
$ws.Range("N126").Value = -10040
